$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2; this shifts the existing rows 2-21 down to rows 3-22
$ws.Rows("2:2").Insert()
$ws.Range("A2:H2").ClearFormats()

# New row 2 sensor reading (ax..gz); timestamp/label set below with the rest
$ws.Range("C2").Value = -0.002831595284598272
$ws.Range("D2").Value = -0.19119707601411
$ws.Range("E2").Value = 0.1513774214046346
$ws.Range("F2").Value = -0.6229298114776611
$ws.Range("G2").Value = -4.338823318481445
$ws.Range("H2").Value = 0.1601994037628173

# 9 brand-new trailing sensor readings appended at rows 23-31
$ws.Range("C23").Value = -2.750307172536849
$ws.Range("D23").Value = -1.468972404088292
$ws.Range("E23").Value = 0.1615269269261984
$ws.Range("F23").Value = 0.5007568001747131
$ws.Range("G23").Value = -1.671174645423889
$ws.Range("H23").Value = -0.2092213481664657

$ws.Range("C24").Value = -3.695782780647288
$ws.Range("D24").Value = -0.6018874943256158
$ws.Range("E24").Value = 1.057165026664715
$ws.Range("F24").Value = 0.0099265603348612
$ws.Range("G24").Value = -3.577379703521729
$ws.Range("H24").Value = -0.4952589869499206

$ws.Range("C25").Value = -5.04273155757359
$ws.Range("D25").Value = 2.776979684829711
$ws.Range("E25").Value = -3.505371774945944
$ws.Range("F25").Value = -1.279915452003479
$ws.Range("G25").Value = -6.508005619049072
$ws.Range("H25").Value = -0.1085812970995903

$ws.Range("C26").Value = -2.32667221341812
$ws.Range("D26").Value = 0.9583009992326832
$ws.Range("E26").Value = -5.241142443248169
$ws.Range("F26").Value = -0.9758572578430176
$ws.Range("G26").Value = -1.915215253829956
$ws.Range("H26").Value = 1.272432327270508

$ws.Range("C27").Value = 2.433021928582876
$ws.Range("D27").Value = -2.759944068534038
$ws.Range("E27").Value = 3.60468020609447
$ws.Range("F27").Value = 0.0806342139840126
$ws.Range("G27").Value = -1.140027284622192
$ws.Range("H27").Value = 0.0485637858510017

$ws.Range("C28").Value = -2.125196490968979
$ws.Range("D28").Value = 0.5359780830996369
$ws.Range("E28").Value = 2.942476987838745
$ws.Range("F28").Value = 0.8413141369819641
$ws.Range("G28").Value = -0.3859141170978546
$ws.Range("H28").Value = 0.5590944290161133

$ws.Range("C29").Value = -3.455752406801498
$ws.Range("D29").Value = 0.2274821900895648
$ws.Range("E29").Value = 2.951905420848299
$ws.Range("F29").Value = -0.1750128865242004
$ws.Range("G29").Value = 2.383749008178711
$ws.Range("H29").Value = -0.1401935666799545

$ws.Range("C30").Value = -2.549259322030204
$ws.Range("D30").Value = 0.480571014540536
$ws.Range("E30").Value = 3.463431903294155
$ws.Range("F30").Value = 0.4100432991981506
$ws.Range("G30").Value = 2.355190992355347
$ws.Range("H30").Value = 0.1059851199388504

$ws.Range("C31").Value = 0.6998523473739624
$ws.Range("D31").Value = -1.159572852775452
$ws.Range("E31").Value = -0.4669593572616656
$ws.Range("F31").Value = -0.6436992287635803
$ws.Range("G31").Value = 0.6840163469314575
$ws.Range("H31").Value = -0.07849618047475811

# timestamp (A) follows a fixed 100ms cadence and label (B) is constant; re-stamp every data row (2-31) since the row insert only shifted columns C:H
for ($r = 2; $r -le 31; $r++) {
    $ws.Cells.Item($r, 1).Value = ($r - 2) * 100
    $ws.Cells.Item($r, 2).Value = "struggle"
}